$wb = $excel.ActiveWorkbook

# ===== Sheet: Demanda =====
$ws = $wb.Worksheets.Item("Hoja1")
$ws.Name = "Demanda"

# -- values --
$ws.Range("A1").Value = "T"
$ws.Range("B1").Value = "R1"
$ws.Range("C1").Value = "R2"
$ws.Range("D1").Value = "R3"
$ws.Range("E1").Value = "R4"
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 1000
$ws.Range("C3").Value = 1500
$ws.Range("D3").Value = 800
$ws.Range("E3").Value = 2000
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 1200
$ws.Range("C4").Value = 1500
$ws.Range("D4").Value = 800
$ws.Range("E4").Value = 2500
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 1350
$ws.Range("C5").Value = 1750
$ws.Range("D5").Value = 1000
$ws.Range("E5").Value = 3000

# -- formatting --
foreach ($a in @("A3","C3","D3","E3","A4","C4","D4","E4")) {
  $c = $ws.Range($a)
  $c.Font.Name = "Times New Roman"
  $c.Font.Size = 11
  $c.VerticalAlignment = -4108
  $c.WrapText = $true
}
foreach ($a in @("A1")) {
  $c = $ws.Range($a)
  $c.Font.Name = "Times New Roman"
  $c.Font.Size = 11
  $c.VerticalAlignment = -4108
  $c.WrapText = $true
  $c.Borders.Item(8).LineStyle = 1
  $c.Borders.Item(8).Weight = -4138
}
foreach ($a in @("B1","C1","D1","E1","A5","C5","D5","E5")) {
  $c = $ws.Range($a)
  $c.Font.Name = "Times New Roman"
  $c.Font.Size = 11
  $c.VerticalAlignment = -4108
  $c.WrapText = $true
  $c.Borders.Item(9).LineStyle = 1
  $c.Borders.Item(9).Weight = -4138
}
foreach ($a in @("B3")) {
  $c = $ws.Range($a)
  $c.Font.Name = "Times New Roman"
  $c.Font.Size = 11
  $c.VerticalAlignment = -4108
  $c.WrapText = $true
  $c.Borders.Item(7).LineStyle = 1
  $c.Borders.Item(7).Weight = 2
  $c.Borders.Item(8).LineStyle = 1
  $c.Borders.Item(8).Weight = -4138
}
foreach ($a in @("B4")) {
  $c = $ws.Range($a)
  $c.Font.Name = "Times New Roman"
  $c.Font.Size = 11
  $c.VerticalAlignment = -4108
  $c.WrapText = $true
  $c.Borders.Item(7).LineStyle = 1
  $c.Borders.Item(7).Weight = 2
}
foreach ($a in @("B5")) {
  $c = $ws.Range($a)
  $c.Font.Name = "Times New Roman"
  $c.Font.Size = 11
  $c.VerticalAlignment = -4108
  $c.WrapText = $true
  $c.Borders.Item(7).LineStyle = 1
  $c.Borders.Item(7).Weight = 2
  $c.Borders.Item(9).LineStyle = 1
  $c.Borders.Item(9).Weight = -4138
}
foreach ($a in @("A2","B2","C2","D2","E2")) {
  $c = $ws.Range($a)
  $c.Font.Name = "Times New Roman"
  $c.Font.Size = 11
  $c.VerticalAlignment = -4108
  $c.WrapText = $true
}

# -- row heights --
$ws.Rows.Item(1).RowHeight = 15.75
$ws.Rows.Item(2).RowHeight = 15.75
$ws.Rows.Item(5).RowHeight = 15.75

$ws.Range("C13").Select() | Out-Null


# ===== Sheet: CostoPF2PE =====
$ws = $wb.Worksheets.Item("Hoja2")
$ws.Name = "CostoPF2PE"

# -- values --
$ws.Range("A1").Value = "PFX"
$ws.Range("B1").Value = "PE1"
$ws.Range("C1").Value = "PE2"
$ws.Range("A2").Value = "PF1"
$ws.Range("B2").Value = 1.8
$ws.Range("C2").Value = 2.2
$ws.Range("A3").Value = "PF2"
$ws.Range("B3").Value = 2.9
$ws.Range("C3").Value = 0.6
$ws.Range("A4").Value = "PF3"
$ws.Range("B4").Value = 1.5
$ws.Range("C4").Value = 2.25

# -- formatting --
foreach ($a in @("A2","A3")) {
  $c = $ws.Range($a)
  $c.Font.Name = "Times New Roman"
  $c.Font.Size = 11
  $c.VerticalAlignment = -4108
  $c.WrapText = $true
}
foreach ($a in @("A1")) {
  $c = $ws.Range($a)
  $c.Font.Name = "Times New Roman"
  $c.Font.Size = 11
  $c.VerticalAlignment = -4108
  $c.WrapText = $true
  $c.Borders.Item(8).LineStyle = 1
  $c.Borders.Item(8).Weight = -4138
}
foreach ($a in @("B1","C1","A4")) {
  $c = $ws.Range($a)
  $c.Font.Name = "Times New Roman"
  $c.Font.Size = 11
  $c.VerticalAlignment = -4108
  $c.WrapText = $true
  $c.Borders.Item(9).LineStyle = 1
  $c.Borders.Item(9).Weight = -4138
}
foreach ($a in @("C2","C3")) {
  $c = $ws.Range($a)
  $c.Font.Name = "Times New Roman"
  $c.Font.Size = 11
  $c.VerticalAlignment = -4108
  $c.WrapText = $true
  $c.NumberFormat = '"$"#,##0.00_);[Red]\("$"#,##0.00\)'
}
foreach ($a in @("C4")) {
  $c = $ws.Range($a)
  $c.Font.Name = "Times New Roman"
  $c.Font.Size = 11
  $c.VerticalAlignment = -4108
  $c.WrapText = $true
  $c.NumberFormat = '"$"#,##0.00_);[Red]\("$"#,##0.00\)'
  $c.Borders.Item(9).LineStyle = 1
  $c.Borders.Item(9).Weight = -4138
}
foreach ($a in @("B2")) {
  $c = $ws.Range($a)
  $c.Font.Name = "Times New Roman"
  $c.Font.Size = 11
  $c.VerticalAlignment = -4108
  $c.WrapText = $true
  $c.NumberFormat = '"$"#,##0.00_);[Red]\("$"#,##0.00\)'
  $c.Borders.Item(7).LineStyle = 1
  $c.Borders.Item(7).Weight = 2
  $c.Borders.Item(8).LineStyle = 1
  $c.Borders.Item(8).Weight = -4138
}
foreach ($a in @("B3")) {
  $c = $ws.Range($a)
  $c.Font.Name = "Times New Roman"
  $c.Font.Size = 11
  $c.VerticalAlignment = -4108
  $c.WrapText = $true
  $c.NumberFormat = '"$"#,##0.00_);[Red]\("$"#,##0.00\)'
  $c.Borders.Item(7).LineStyle = 1
  $c.Borders.Item(7).Weight = 2
}
foreach ($a in @("B4")) {
  $c = $ws.Range($a)
  $c.Font.Name = "Times New Roman"
  $c.Font.Size = 11
  $c.VerticalAlignment = -4108
  $c.WrapText = $true
  $c.NumberFormat = '"$"#,##0.00_);[Red]\("$"#,##0.00\)'
  $c.Borders.Item(7).LineStyle = 1
  $c.Borders.Item(7).Weight = 2
  $c.Borders.Item(9).LineStyle = 1
  $c.Borders.Item(9).Weight = -4138
}

# -- row heights --
$ws.Rows.Item(1).RowHeight = 15.75
$ws.Rows.Item(4).RowHeight = 15.75

$ws.Range("A25").Select() | Out-Null


# ===== Sheet: CostoPE2R =====
$ws = $wb.Worksheets.Item("Hoja3")
$ws.Name = "CostoPE2R"

# -- values --
$ws.Range("A1").Value = "PEX"
$ws.Range("B1").Value = "R1"
$ws.Range("C1").Value = "R2"
$ws.Range("D1").Value = "R3"
$ws.Range("E1").Value = "R4"
$ws.Range("A2").Value = "PE1"
$ws.Range("B2").Value = 1.1
$ws.Range("C2").Value = 2.15
$ws.Range("D2").Value = 5.4
$ws.Range("E2").Value = 6.1
$ws.Range("A3").Value = "PE2"
$ws.Range("B3").Value = 7.15
$ws.Range("C3").Value = 6.7
$ws.Range("D3").Value = 4.5
$ws.Range("E3").Value = 2.15

# -- formatting --
foreach ($a in @("A1")) {
  $c = $ws.Range($a)
  $c.Font.Name = "Times New Roman"
  $c.Font.Size = 11
  $c.VerticalAlignment = -4108
  $c.WrapText = $true
  $c.Borders.Item(8).LineStyle = 1
  $c.Borders.Item(8).Weight = -4138
}
foreach ($a in @("B1","C1","D1","E1","A2","A3")) {
  $c = $ws.Range($a)
  $c.Font.Name = "Times New Roman"
  $c.Font.Size = 11
  $c.VerticalAlignment = -4108
  $c.WrapText = $true
  $c.Borders.Item(9).LineStyle = 1
  $c.Borders.Item(9).Weight = -4138
}
foreach ($a in @("C2","D2","E2","C3","D3","E3")) {
  $c = $ws.Range($a)
  $c.Font.Name = "Times New Roman"
  $c.Font.Size = 11
  $c.VerticalAlignment = -4108
  $c.WrapText = $true
  $c.NumberFormat = '"$"#,##0.00_);[Red]\("$"#,##0.00\)'
  $c.Borders.Item(9).LineStyle = 1
  $c.Borders.Item(9).Weight = -4138
}
foreach ($a in @("B3")) {
  $c = $ws.Range($a)
  $c.Font.Name = "Times New Roman"
  $c.Font.Size = 11
  $c.VerticalAlignment = -4108
  $c.WrapText = $true
  $c.NumberFormat = '"$"#,##0.00_);[Red]\("$"#,##0.00\)'
  $c.Borders.Item(7).LineStyle = 1
  $c.Borders.Item(7).Weight = 2
  $c.Borders.Item(9).LineStyle = 1
  $c.Borders.Item(9).Weight = -4138
}
foreach ($a in @("B2")) {
  $c = $ws.Range($a)
  $c.Font.Name = "Times New Roman"
  $c.Font.Size = 11
  $c.VerticalAlignment = -4108
  $c.WrapText = $true
  $c.NumberFormat = '"$"#,##0.00_);[Red]\("$"#,##0.00\)'
  $c.Borders.Item(7).LineStyle = 1
  $c.Borders.Item(7).Weight = 2
  $c.Borders.Item(8).LineStyle = 1
  $c.Borders.Item(8).Weight = -4138
  $c.Borders.Item(9).LineStyle = 1
  $c.Borders.Item(9).Weight = -4138
}

# -- row heights --
$ws.Rows.Item(1).RowHeight = 15.75
$ws.Rows.Item(2).RowHeight = 15.75
$ws.Rows.Item(3).RowHeight = 15.75

# -- column widths --
$ws.Columns.Item(1).ColumnWidth = 10.42578125

$ws.Range("E3").Select() | Out-Null


# -- restore active sheet/tab --
$ws = $wb.Worksheets.Item("Demanda")
$ws.Activate()
$ws.Range("C13").Select() | Out-Null
